# Generate Report for Handback
#
# The CI report regenerated two rows worth of handback metadata: the first
# row's source markdown file got a new GUID (1054cb0d... -> b97a8a6b...) and
# the second row's source markdown got swapped for an entirely new file
# (de78cca7... -> ffff65923d38...). New xliff file names and new timestamps
# were produced by the (re)generation run as well.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Sheets.Item("Overview")
$wsZhCn     = $wb.Sheets.Item("zh-cn")
$wsDeDe     = $wb.Sheets.Item("de-de")

# ---- Overview sheet ----------------------------------------------------
$wsOverview.Range("A2").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.md"
$wsOverview.Range("B2").Value = "e2e\b97a8a6b-2738-43ca-9338-babba44e2b94.md"
$wsOverview.Range("G2").Value = "2016-08-29 15:13:23"

$wsOverview.Range("A3").Value = "ffff65923d38-c131-4f9d-b869-90ae172a464d.md"
$wsOverview.Range("B3").Value = "e2e\ffff65923d38-c131-4f9d-b869-90ae172a464d.md"
$wsOverview.Range("G3").Value = "2016-08-29 15:13:23"

# ---- zh-cn sheet ---------------------------------------------------------
$wsZhCn.Range("A2").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.md"
$wsZhCn.Range("G2").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.79f01d5ff72c355ab3e762393336a1512850592a.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-29 15:13:18"
$wsZhCn.Range("I2").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.md"
$wsZhCn.Range("J2").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.79f01d5ff72c355ab3e762393336a1512850592a.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-29 15:13:35"

$wsZhCn.Range("A3").Value = "ffff65923d38-c131-4f9d-b869-90ae172a464d.md"
$wsZhCn.Range("G3").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.79f01d5ff72c355ab3e762393336a1512850592a.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-08-29 15:13:18"
$wsZhCn.Range("I3").Value = "ffff65923d38-c131-4f9d-b869-90ae172a464d.md"
$wsZhCn.Range("J3").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.79f01d5ff72c355ab3e762393336a1512850592a.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-08-29 15:13:35"

# ---- de-de sheet ---------------------------------------------------------
$wsDeDe.Range("A2").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.md"
$wsDeDe.Range("G2").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.79f01d5ff72c355ab3e762393336a1512850592a.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-29 15:13:23"
$wsDeDe.Range("I2").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.md"
$wsDeDe.Range("J2").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.79f01d5ff72c355ab3e762393336a1512850592a.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-29 15:13:42"

$wsDeDe.Range("A3").Value = "ffff65923d38-c131-4f9d-b869-90ae172a464d.md"
$wsDeDe.Range("G3").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.79f01d5ff72c355ab3e762393336a1512850592a.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-08-29 15:13:23"
$wsDeDe.Range("I3").Value = "ffff65923d38-c131-4f9d-b869-90ae172a464d.md"
$wsDeDe.Range("J3").Value = "b97a8a6b-2738-43ca-9338-babba44e2b94.79f01d5ff72c355ab3e762393336a1512850592a.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-08-29 15:13:42"

# ---- Refresh hyperlink display text (targets/rIds are left untouched, --
# ---- matching the real handback tool which re-writes the report table --
# ---- without re-resolving the github blob URLs) -------------------------

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31d81c69490e4603d6012daa49fb182d6719a490/e2e/1054cb0d-4534-4036-bb4b-ea26a3c7dfc9.md", "", "", "e2e\b97a8a6b-2738-43ca-9338-babba44e2b94.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31d81c69490e4603d6012daa49fb182d6719a490/e2e/de78cca7-9532-4156-a74c-fdd2e9eea052.md", "", "", "e2e\ffff65923d38-c131-4f9d-b869-90ae172a464d.md")

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31d81c69490e4603d6012daa49fb182d6719a490/e2e/1054cb0d-4534-4036-bb4b-ea26a3c7dfc9.md", "", "", "b97a8a6b-2738-43ca-9338-babba44e2b94.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e33cd45215d90a20809654119808020616ac9f67/e2e/1054cb0d-4534-4036-bb4b-ea26a3c7dfc9.md", "", "", "b97a8a6b-2738-43ca-9338-babba44e2b94.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31d81c69490e4603d6012daa49fb182d6719a490/e2e/de78cca7-9532-4156-a74c-fdd2e9eea052.md", "", "", "ffff65923d38-c131-4f9d-b869-90ae172a464d.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e33cd45215d90a20809654119808020616ac9f67/e2e/de78cca7-9532-4156-a74c-fdd2e9eea052.md", "", "", "ffff65923d38-c131-4f9d-b869-90ae172a464d.md")

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31d81c69490e4603d6012daa49fb182d6719a490/e2e/1054cb0d-4534-4036-bb4b-ea26a3c7dfc9.md", "", "", "b97a8a6b-2738-43ca-9338-babba44e2b94.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/dcfda62f189468e9a1400b52c6376f60a5cfd1cd/e2e/1054cb0d-4534-4036-bb4b-ea26a3c7dfc9.md", "", "", "b97a8a6b-2738-43ca-9338-babba44e2b94.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/31d81c69490e4603d6012daa49fb182d6719a490/e2e/de78cca7-9532-4156-a74c-fdd2e9eea052.md", "", "", "ffff65923d38-c131-4f9d-b869-90ae172a464d.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/dcfda62f189468e9a1400b52c6376f60a5cfd1cd/e2e/de78cca7-9532-4156-a74c-fdd2e9eea052.md", "", "", "ffff65923d38-c131-4f9d-b869-90ae172a464d.md")
